# "Generate Report for Handoff"
# Updates the localization-status report to reflect that the content has
# moved from "In Translation" to "Ready for handoff", refreshes the
# handoff/generation timestamps, and widens the datetime columns so the
# new (longer) status text / timestamps are fully visible.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet -------------------------------------------------------
# zh-cn / de-de status columns
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
# Latest HO Xliff Generate Date
$wsOverview.Range("G2").Value = "2016-08-16 02:35:58"

# --- zh-cn sheet ------------------------------------------------------------
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-16 02:35:53"

# --- de-de sheet ------------------------------------------------------------
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-16 02:35:58"

# --- Widen the datetime columns to fit the refreshed values -----------------
# Target stored width is 17.2159881591797 characters; this engine stores
# column widths as whole pixels (6 px/char at the default font), so the
# nearest representable ColumnWidth is 98 px -> 98/6 characters.
$newDateColWidth = 98 / 6

$wsOverview.Columns.Item(5).ColumnWidth = $newDateColWidth  # column E
$wsOverview.Columns.Item(6).ColumnWidth = $newDateColWidth  # column F

$wsZhCn.Columns.Item(3).ColumnWidth = $newDateColWidth      # column C

$wsDeDe.Columns.Item(3).ColumnWidth = $newDateColWidth      # column C
